# Updated cryptos list on Thu Nov 30 22:10:21 UTC 2023 with GitHub Actions
# Refresh the scraped Price (D) and Volume(1h) (E) columns for every coin row.
#
# D-column values look numeric ("5.90", "227.63", ...) but must stay plain
# text (matching the sheet's original inlineStr cells), otherwise Excel's
# COM layer would silently coerce them to Number and drop the formatting
# (e.g. "5.90" -> 5.9). Force text via NumberFormat "@" just for the write,
# then restore the cell's original Style so no stray formatting is left
# behind.
function Set-TextValue($range, [string]$value) {
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $savedStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Set-TextValue $ws.Range("D2") "37.784.51"
$ws.Range("E2").Value = "  -0.09%  "
Set-TextValue $ws.Range("D3") "2.047.32"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  +0.10%  "
Set-TextValue $ws.Range("D5") "227.63"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("E6").Value = "  -0.39%  "
Set-TextValue $ws.Range("D7") "59.88"
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("E8").Value = "  +0.05%  "
Set-TextValue $ws.Range("D9") "0.377"
$ws.Range("E9").Value = "  -1.61%  "
Set-TextValue $ws.Range("D10") "0.0839"
$ws.Range("E10").Value = "  +3.24%  "
$ws.Range("E11").Value = "  +0.21%  "
Set-TextValue $ws.Range("D12") "2.350.09"
$ws.Range("E12").Value = "  +0.82%  "
Set-TextValue $ws.Range("D13") "14.39"
$ws.Range("E13").Value = "  -1.29%  "
Set-TextValue $ws.Range("D14") "21.36"
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("E15").Value = "  +6.26%  "
$ws.Range("E16").Value = "  +0.77%  "
Set-TextValue $ws.Range("D17") "2.038.55"
$ws.Range("E17").Value = "  +0.70%  "
Set-TextValue $ws.Range("D18") "37.758.82"
$ws.Range("E18").Value = "  +0.18%  "
Set-TextValue $ws.Range("D19") "69.42"
$ws.Range("E19").Value = "  -0.66%  "
Set-TextValue $ws.Range("D20") "5.90"
$ws.Range("E20").Value = "  -2.00%  "
Set-TextValue $ws.Range("D21") "0.0₃0830"
$ws.Range("E21").Value = "  +0.66%  "
Set-TextValue $ws.Range("D22") "222.41"
$ws.Range("E22").Value = "  -1.21%  "
$ws.Range("E23").Value = "  +0.46%  "
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("E25").Value = "  +3.14%  "
Set-TextValue $ws.Range("D26") "168.87"
$ws.Range("E26").Value = "  +2.24%  "
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("E28").Value = "  -0.37%  "
Set-TextValue $ws.Range("D29") "18.77"
$ws.Range("E29").Value = "  -0.95%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("E32").Value = "  +8.20%  "
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("E34").Value = "  +0.62%  "
Set-TextValue $ws.Range("D35") "0.0602"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  +2.09%  "
$ws.Range("E37").Value = "  +4.39%  "
Set-TextValue $ws.Range("D38") "3.47"
$ws.Range("E38").Value = "  +7.30%  "
$ws.Range("E39").Value = "  -0.03%  "
Set-TextValue $ws.Range("D40") "18.33"
$ws.Range("E40").Value = "  +9.25%  "
Set-TextValue $ws.Range("D41") "1.528.62"
$ws.Range("E41").Value = "  -0.07%  "
Set-TextValue $ws.Range("D42") "97.71"
$ws.Range("E42").Value = "  +0.87%  "
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("E45").Value = "  +0.47%  "
Set-TextValue $ws.Range("D46") "0.0890"
$ws.Range("E46").Value = "  -2.91%  "
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("E48").Value = "  -0.10%  "
Set-TextValue $ws.Range("D49") "2.95"
$ws.Range("E49").Value = "  -0.39%  "
Set-TextValue $ws.Range("D50") "7.10"
$ws.Range("E50").Value = "  +0.02%  "
Set-TextValue $ws.Range("D51") "2.238.56"
$ws.Range("E51").Value = "  +0.85%  "
